$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the opening H1
#    heading ("Play Age of Halvar Free: Unique Wild Halvar Feature").
# ---------------------------------------------------------------------------
$heading = $d.Paragraphs.Item(1)
$heading.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"
$metaStart = $metaPara.Range.Start

$boldLabel = "Meta description"
$restOfText = ": Read a review of Age of Halvar, an online slot game with a unique Wild Halvar feature. Play for free and enjoy the fun Nordic setting and atmosphere."

# Type the full sentence first (plain run) ...
$textRange = $d.Range($metaStart, $metaStart)
$textRange.InsertAfter($boldLabel + $restOfText)

# ... then bold just the "Meta description" label, splitting it into its own run.
$labelRange = $d.Range($metaStart, $metaStart + $boldLabel.Length)
$labelRange.Bold = 1

# Finally, prepend the empty leading run that matches the document's usual
# paragraph shape (<w:r/> followed by the real runs).
$leadStart = $metaPara.Range.Start
$emptyRunXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$leadRange = $d.Range($leadStart, $leadStart)
$leadRange.InsertXML($emptyRunXml)

# ---------------------------------------------------------------------------
# 2) Drop the now-duplicated bold "Play Age of Halvar Free..." paragraph that
#    used to sit just before the closing italic paragraph. (The opening H1
#    heading has the same text, so disambiguate on paragraph style: the
#    duplicate to remove is a plain "Normal" paragraph, not "Heading 1".)
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Normal" -and $p.Range.Text.Trim() -eq "Play Age of Halvar Free: Unique Wild Halvar Feature") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Replace the closing italic paragraph's text with the new (garbled)
#    image-prompt content, keeping its italic formatting intact. Target it
#    directly (it's the document's last paragraph) so the still-identical
#    sentence inside the new Meta description paragraph is left untouched.
# ---------------------------------------------------------------------------
$newClosing = 'Create a feature image fitting the game "Age Of Halvar": - The image should be in cartoon style - The should feature a happy Maya warrior with glasses Sorry, there seems to be a confusion in your prompt. The game is called "Age Of Halvar", which is based on Vikings theme, but your prompt is asking for a feature image of a happy Maya warrior with glasses. Please clarify the prompt so I can provide an appropriate response.'

$closingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$closingRange = $d.Range($closingPara.Range.Start, $closingPara.Range.End - 1)
$closingRange.Text = $newClosing

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
